$d = $word.ActiveDocument

# Delete paragraphs 2, 3, and 4 (the duplicated/extra content added by the
# "Added More Stuff to MyNewBook" commit), restoring the document to just
# the original first paragraph followed by the trailing empty paragraph.
for ($i = 4; $i -ge 2; $i--) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Delete()
}
